$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 3583  # was 2874.5
$ws.Range("H62").Value = 9466.666999999999  # was 9700.75
$ws.Range("I62").Value = 9466.666999999999  # was 9700.75
$ws.Range("K62").Value = 9466.666999999999  # was 9700.75
$ws.Range("M62").Value = -8842.666999999999  # was -9076.75
$ws.Range("H65").Value = 9466.666999999999  # was 9700.75
$ws.Range("I65").Value = 9466.666999999999  # was 9700.75
$ws.Range("K65").Value = 47333.335  # was 48503.75
$ws.Range("M65").Value = -44213.335  # was -45383.75
$ws.Range("H70").Value = 5592.4243  # was 5732.8125
$ws.Range("J70").Value = 6807.2856  # was 7092.65
$ws.Range("L70").Value = 20421.8568  # was 21277.95
$ws.Range("N70").Value = -20961.8568  # was -21817.95
$ws.Range("H73").Value = 5592.4243  # was 5732.8125
$ws.Range("J73").Value = 6807.2856  # was 7092.65
$ws.Range("L73").Value = 20421.8568  # was 21277.95
$ws.Range("N73").Value = -22293.8568  # was -23149.95
$ws.Range("H82").Value = 719.3333  # was 761.5714
$ws.Range("I82").Value = 586.8  # was 589.8
$ws.Range("J82").Value = 1382  # was 1191
$ws.Range("K82").Value = 1760.4  # was 1769.4
$ws.Range("L82").Value = 4146  # was 3573
$ws.Range("M82").Value = -1354.4  # was -1363.4
$ws.Range("N82").Value = -4958  # was -4385
$ws.Range("H85").Value = 719.3333  # was 761.5714
$ws.Range("I85").Value = 586.8  # was 589.8
$ws.Range("J85").Value = 1382  # was 1191
$ws.Range("K85").Value = 1760.4  # was 1769.4
$ws.Range("L85").Value = 4146  # was 3573
$ws.Range("M85").Value = -356.3999999999999  # was -365.3999999999999
$ws.Range("N85").Value = -6954  # was -6381
$ws.Range("H86").Value = 3003  # was 3002.75
$ws.Range("I86").Value = 3003  # was 3002.75
$ws.Range("K86").Value = 3003  # was 3002.75
$ws.Range("M86").Value = -1880  # was -1879.75
$ws.Range("H89").Value = 3003  # was 3002.75
$ws.Range("I89").Value = 3003  # was 3002.75
$ws.Range("K89").Value = 15015  # was 15013.75
$ws.Range("M89").Value = -9399  # was -9397.75
$ws.Range("H131").Value = 1472  # was 1597.5
$ws.Range("I131").Value = 1472  # was 1597.5
$ws.Range("K131").Value = 4416  # was 4792.5
$ws.Range("M131").Value = 624  # was 247.5
$ws.Range("H132").Value = 2040.625  # was 2048.9375
$ws.Range("I132").Value = 2058.5957  # was 2086.0217
$ws.Range("K132").Value = 6175.7871  # was 6258.0651
$ws.Range("M132").Value = -3645.7871  # was -3728.0651
$ws.Range("H137").Value = 241736.39  # was 241814.61
$ws.Range("I137").Value = 2547.8  # was 2657.3
$ws.Range("K137").Value = 7643.400000000001  # was 7971.900000000001
$ws.Range("M137").Value = -5093.400000000001  # was -5421.900000000001
$ws.Range("H138").Value = 2549.422  # was 2260.7036
$ws.Range("I138").Value = 917.96  # was 791
$ws.Range("J138").Value = 4588.75  # was 4570.2383
$ws.Range("K138").Value = 2753.88  # was 2373
$ws.Range("L138").Value = 13766.25  # was 13710.7149
$ws.Range("M138").Value = 2386.12  # was 2767
$ws.Range("N138").Value = -24046.25  # was -23990.7149
$ws.Range("H141").Value = 704.4054  # was 718.4722
$ws.Range("I141").Value = 723.5454999999999  # was 739.96875
$ws.Range("K141").Value = 2170.6365  # was 2219.90625
$ws.Range("M141").Value = 3009.3635  # was 2960.09375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16643.654  # was 16784.69
$ws.Range("I32").Value = 13170  # was 14499.5
$ws.Range("J32").Value = 24362.889  # was 20524.092
$ws.Range("K32").Value = 13170  # was 14499.5
$ws.Range("L32").Value = 24362.889  # was 20524.092
$ws.Range("M32").Value = -12883  # was -14212.5
$ws.Range("N32").Value = -24936.889  # was -21098.092
$ws.Range("H61").Value = 17505.5  # was 16056.442
$ws.Range("I61").Value = 2446.5454  # was 2200.9216
$ws.Range("J61").Value = 50635.2  # was 53247.58
$ws.Range("K61").Value = 2446.5454  # was 2200.9216
$ws.Range("L61").Value = 50635.2  # was 53247.58
$ws.Range("M61").Value = -2234.5454  # was -1988.9216
$ws.Range("N61").Value = -51059.2  # was -53671.58
$ws.Range("H74").Value = 83908.27  # was 48793.105
$ws.Range("I74").Value = 54898.176  # was 29402.844
$ws.Range("J74").Value = 182542.6  # was 152207.83
$ws.Range("K74").Value = 54898.176  # was 29402.844
$ws.Range("L74").Value = 182542.6  # was 152207.83
$ws.Range("M74").Value = -54024.176  # was -28528.844
$ws.Range("N74").Value = -184290.6  # was -153955.83
$ws.Range("H77").Value = 83908.27  # was 48793.105
$ws.Range("I77").Value = 54898.176  # was 29402.844
$ws.Range("J77").Value = 182542.6  # was 152207.83
$ws.Range("K77").Value = 274490.88  # was 147014.22
$ws.Range("L77").Value = 912713  # was 761039.1499999999
$ws.Range("M77").Value = -270122.88  # was -142646.22
$ws.Range("N77").Value = -921449  # was -769775.1499999999
$ws.Range("H106").Value = 0  # was 50185
$ws.Range("J106").Value = 0  # was 50185
$ws.Range("L106").Value = 0  # was 50185
$ws.Range("N106").ClearContents()  # was -52709
$ws.Range("H122").Value = 62148.133  # was 41006.566
$ws.Range("I122").Value = 1851.8334  # was 1691.8948
$ws.Range("J122").Value = 303333.34  # was 227751.25
$ws.Range("K122").Value = 5555.5002  # was 5075.6844
$ws.Range("L122").Value = 910000.02  # was 683253.75
$ws.Range("M122").Value = -3105.5002  # was -2625.6844
$ws.Range("N122").Value = -914900.02  # was -688153.75
$ws.Range("H132").Value = 2903.5083  # was 2844.8032
$ws.Range("I132").Value = 2578.2363  # was 2536.0715
$ws.Range("J132").Value = 5885.1665  # was 6302.6
$ws.Range("K132").Value = 7734.7089  # was 7608.2145
$ws.Range("L132").Value = 17655.4995  # was 18907.8
$ws.Range("M132").Value = -5204.7089  # was -5078.2145
$ws.Range("N132").Value = -22715.4995  # was -23967.8
$ws.Range("H136").Value = 17505.5  # was 16056.442
$ws.Range("I136").Value = 2446.5454  # was 2200.9216
$ws.Range("J136").Value = 50635.2  # was 53247.58
$ws.Range("K136").Value = 7339.6362  # was 6602.764800000001
$ws.Range("L136").Value = 151905.6  # was 159742.74
$ws.Range("M136").Value = -4789.6362  # was -4052.764800000001
$ws.Range("N136").Value = -157005.6  # was -164842.74

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4393.5127  # was 3276.7068
$ws.Range("I134").Value = 3952.76  # was 2670.932
$ws.Range("K134").Value = 11858.28  # was 8012.795999999999
$ws.Range("M134").Value = -9323.280000000001  # was -5477.795999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2851.4478  # was 2862.8484
$ws.Range("I31").Value = 2768.0435  # was 2794.2856
$ws.Range("J31").Value = 2895.0454  # was 2894.8445
$ws.Range("K31").Value = 2768.0435  # was 2794.2856
$ws.Range("L31").Value = 2895.0454  # was 2894.8445
$ws.Range("M31").Value = -2473.0435  # was -2499.2856
$ws.Range("N31").Value = -3485.0454  # was -3484.8445
$ws.Range("H34").Value = 2851.4478  # was 2862.8484
$ws.Range("I34").Value = 2768.0435  # was 2794.2856
$ws.Range("J34").Value = 2895.0454  # was 2894.8445
$ws.Range("K34").Value = 2768.0435  # was 2794.2856
$ws.Range("L34").Value = 2895.0454  # was 2894.8445
$ws.Range("M34").Value = -2566.0435  # was -2592.2856
$ws.Range("N34").Value = -3299.0454  # was -3298.8445
$ws.Range("H88").Value = 29518.5  # was 30351.834
$ws.Range("J88").Value = 28960  # was 29960
$ws.Range("L88").Value = 28960  # was 29960
$ws.Range("N88").Value = -29772  # was -30772
$ws.Range("H91").Value = 29518.5  # was 30351.834
$ws.Range("J91").Value = 28960  # was 29960
$ws.Range("L91").Value = 28960  # was 29960
$ws.Range("N91").Value = -31768  # was -32768
$ws.Range("H107").Value = 38290.2  # was 43470.91
$ws.Range("I107").Value = 48505.156  # was 57543.938
$ws.Range("K107").Value = 48505.156  # was 57543.938
$ws.Range("M107").Value = -46585.156  # was -55623.938
$ws.Range("H132").Value = 5111.72  # was 5524.8696
$ws.Range("I132").Value = 1628.7  # was 1769.6111
$ws.Range("K132").Value = 4886.1  # was 5308.8333
$ws.Range("M132").Value = -2356.1  # was -2778.8333
$ws.Range("H134").Value = 2859.653  # was 2917.1875
$ws.Range("I134").Value = 2592.348  # was 2647.7778
$ws.Range("K134").Value = 7777.044  # was 7943.3334
$ws.Range("M134").Value = -5242.044  # was -5408.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2656  # was 787.3043
$ws.Range("I12").Value = 0  # was 134.5
$ws.Range("J12").Value = 2656  # was 1135.4667
$ws.Range("K12").Value = 0  # was 403.5
$ws.Range("L12").Value = 7968  # was 3406.4001
$ws.Range("M12").ClearContents()  # was -230.5
$ws.Range("N12").Value = -8314  # was -3752.4001
$ws.Range("H13").Value = 2197.6  # was 2034.1818
$ws.Range("I13").Value = 1036.2  # was 930.1667
$ws.Range("K13").Value = 3108.6  # was 2790.5001
$ws.Range("M13").Value = -2940.6  # was -2622.5001
$ws.Range("H113").Value = 1475.7059  # was 1426.8334
$ws.Range("I113").Value = 785  # was 761.375
$ws.Range("K113").Value = 2355  # was 2284.125
$ws.Range("M113").Value = -185  # was -114.125
$ws.Range("H140").Value = 1613  # was 1597.675
$ws.Range("I140").Value = 1341.5428  # was 1332.0555
$ws.Range("K140").Value = 4024.6284  # was 3996.1665
$ws.Range("M140").Value = 1155.3716  # was 1183.8335

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10334.8  # was 8962.333000000001
$ws.Range("I80").Value = 3149.7778  # was 3054.8
$ws.Range("J80").Value = 75000  # was 38500
$ws.Range("K80").Value = 3149.7778  # was 3054.8
$ws.Range("L80").Value = 75000  # was 38500
$ws.Range("M80").Value = -2151.7778  # was -2056.8
$ws.Range("N80").Value = -76996  # was -40496
$ws.Range("H83").Value = 10334.8  # was 8962.333000000001
$ws.Range("I83").Value = 3149.7778  # was 3054.8
$ws.Range("J83").Value = 75000  # was 38500
$ws.Range("K83").Value = 15748.889  # was 15274
$ws.Range("L83").Value = 375000  # was 192500
$ws.Range("M83").Value = -10756.889  # was -10282
$ws.Range("N83").Value = -384984  # was -202484
$ws.Range("H132").Value = 3747.3096  # was 3503.9348
$ws.Range("I132").Value = 3858.8965  # was 3505.8484
$ws.Range("J132").Value = 3498.3845  # was 3499.077
$ws.Range("K132").Value = 11576.6895  # was 10517.5452
$ws.Range("L132").Value = 10495.1535  # was 10497.231
$ws.Range("M132").Value = -9046.6895  # was -7987.5452
$ws.Range("N132").Value = -15555.1535  # was -15557.231

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 41665.875  # was 99999
$ws.Range("J60").Value = 41665.875  # was 99999
$ws.Range("L60").Value = 41665.875  # was 99999
$ws.Range("N60").Value = -42683.875  # was -101017
$ws.Range("H136").Value = 23943.762  # was 23978.096
$ws.Range("I136").Value = 2251.5151  # was 2295.2122
$ws.Range("K136").Value = 6754.5453  # was 6885.6366
$ws.Range("M136").Value = -4204.5453  # was -4335.6366

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 36848.5  # was 34478.8
$ws.Range("J80").Value = 36848.5  # was 34478.8
$ws.Range("L80").Value = 36848.5  # was 34478.8
$ws.Range("N80").Value = -38844.5  # was -36474.8
$ws.Range("H83").Value = 36848.5  # was 34478.8
$ws.Range("J83").Value = 36848.5  # was 34478.8
$ws.Range("L83").Value = 110545.5  # was 103436.4
$ws.Range("N83").Value = -120529.5  # was -113420.4
$ws.Range("H122").Value = 3825.2046  # was 3895.6223
$ws.Range("I122").Value = 3759.6099  # was 3771.6584
$ws.Range("J122").Value = 4721.6665  # was 5166.25
$ws.Range("K122").Value = 11278.8297  # was 11314.9752
$ws.Range("L122").Value = 14164.9995  # was 15498.75
$ws.Range("M122").Value = -8828.8297  # was -8864.975199999999
$ws.Range("N122").Value = -19064.9995  # was -20398.75
$ws.Range("H132").Value = 2249.195  # was 2359.8206
$ws.Range("I132").Value = 2205.35  # was 2316.5789
$ws.Range("K132").Value = 6616.049999999999  # was 6949.736699999999
$ws.Range("M132").Value = -4086.049999999999  # was -4419.736699999999
$ws.Range("H136").Value = 3155.87  # was 3220.8572
$ws.Range("I136").Value = 3053.1475  # was 3135.1475
$ws.Range("J136").Value = 3547.5  # was 3547.625
$ws.Range("K136").Value = 9159.442500000001  # was 9405.442500000001
$ws.Range("L136").Value = 10642.5  # was 10642.875
$ws.Range("M136").Value = -6609.442500000001  # was -6855.442500000001
$ws.Range("N136").Value = -15742.5  # was -15742.875

Write-Host "Applied all Phoenix_Profits market-data updates"